$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 22: 2025-12-11 (46002), 四方坪站
$ws.Range("A22").Value = 46002
$ws.Range("B22").Value = "四方坪站"
$ws.Range("C22").Value = 8172.82
$ws.Range("D22").Value = 7014.37
$ws.Range("E22").Value = 2728.39
$ws.Range("F22").Value = 372

# Add new row 23: 2025-12-11 (46002), 高岭站
$ws.Range("A23").Value = 46002
$ws.Range("B23").Value = "高岭站"
$ws.Range("C23").Value = 4169.42
$ws.Range("D23").Value = 3509.36
$ws.Range("E23").Value = 1137.06
$ws.Range("F23").Value = 159

# Update the active cell selection to match where the user ended up (I22)
[void]$ws.Range("I22").Select()
